$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 68: "Create HelloWorld Function" (TypeScript)
$ws.Range("A68").Value = 2667
$ws.Range("B68").Value = "Create HelloWorld Function"
$ws.Range("C68").Value = "TypeScript"
$ws.Range("D68").Value = 'return "Hello World"'

# Row 66: add new Type/Steps cells for the existing "Largest 3-Same-Digit Number in a String" entry
$ws.Range("C66").Value = "List/Loop"
$ws.Range("D66").Value = "C# Contains"

# Update selection to match the new active cell
$ws.Range("D60").Select()
